$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C) for rows 2-6 from 2023-09-06 (45175)
# to 2023-09-14 (45183), preserving existing formatting.
$ws.Range("C2:C6").Value = 45183
